$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scale - EEWW")
Write-Host "Sheet1 name:" $ws.Name
